$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update cell C5 (Max value of pie_threshold_range) from 25 to 20
$ws.Range("C5").Value = 20
